$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.155.69"
Set-TextValue "E2" "  +0.22%  "
Set-TextValue "D3" "1.905.23"
Set-TextValue "E3" "  +0.67%  "
Set-TextValue "D5" "306.69"
Set-TextValue "E5" "  -0.06%  "
Set-TextValue "D6" "1.001"
Set-TextValue "E6" "  +0.19%  "
Set-TextValue "D7" "0.5231"
Set-TextValue "E7" "  +1.65%  "
Set-TextValue "D8" "0.3763"
Set-TextValue "E8" "  +0.22%  "
Set-TextValue "D9" "0.07261"
Set-TextValue "E9" "  +0.75%  "
Set-TextValue "E10" "  -0.13%  "
Set-TextValue "D11" "0.9054"
Set-TextValue "E11" "  +0.02%  "
Set-TextValue "E12" "  +11.04%  "
Set-TextValue "D13" "1.918.94"
Set-TextValue "E13" "  +1.36%  "
Set-TextValue "D14" "96.88"
Set-TextValue "E14" "  +1.95%  "
Set-TextValue "D15" "5.292"
Set-TextValue "E15" "  +0.51%  "
Set-TextValue "D16" "1.001"
Set-TextValue "E16" "  +0.24%  "
Set-TextValue "D17" "0.000008699"
Set-TextValue "E18" "  +0.69%  "
Set-TextValue "D19" "1.001"
Set-TextValue "E19" "  +0.21%  "
Set-TextValue "D20" "27.199.71"
Set-TextValue "E21" "  +0.36%  "
Set-TextValue "D22" "2.157.62"
Set-TextValue "E23" "  +0.66%  "
Set-TextValue "D24" "6.439"
Set-TextValue "E24" "  +0.45%  "
Set-TextValue "D25" "2.307"
Set-TextValue "E25" "  +0.52%  "
Set-TextValue "D26" "146.98"
Set-TextValue "E26" "  +0.91%  "
Set-TextValue "B27" "Toncoin"
Set-TextValue "C27" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D27" "1.753"
Set-TextValue "E27" "  -0.84%  "
Set-TextValue "B28" "EthereumClassic"
Set-TextValue "C28" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D28" "18.26"
Set-TextValue "E28" "  +1.10%  "
Set-TextValue "D29" "115.13"
Set-TextValue "E29" "  +0.52%  "
Set-TextValue "D30" "4.821"
Set-TextValue "E30" "  -0.13%  "
Set-TextValue "D31" "4.917"
Set-TextValue "E31" "  -0.70%  "
Set-TextValue "D32" "0.09315"
Set-TextValue "E32" "  +1.40%  "
Set-TextValue "D33" "0.8001"
Set-TextValue "E33" "  +2.56%  "
Set-TextValue "D34" "0.05057"
Set-TextValue "E34" "  -0.54%  "
Set-TextValue "D35" "1.241"
Set-TextValue "E35" "  +0.31%  "
Set-TextValue "D36" "3.436"
Set-TextValue "E36" "  +4.41%  "
Set-TextValue "D37" "2.948"
Set-TextValue "E37" "  -1.05%  "
Set-TextValue "D38" "2.598"
Set-TextValue "E38" "  -0.79%  "
Set-TextValue "D39" "0.5703"
Set-TextValue "E39" "  +1.87%  "
Set-TextValue "D40" "0.02002"
Set-TextValue "E40" "  +0.17%  "
Set-TextValue "D41" "1.074"
Set-TextValue "E41" "  -0.25%  "
Set-TextValue "D42" "9.103"
Set-TextValue "E42" "  +0.02%  "
Set-TextValue "D43" "6.628"
Set-TextValue "E43" "  -0.60%  "
Set-TextValue "D44" "115.85"
Set-TextValue "E44" "  -1.69%  "
Set-TextValue "D45" "0.1513"
Set-TextValue "E45" "  +0.25%  "
Set-TextValue "D46" "0.4860"
Set-TextValue "E46" "  +0.93%  "
Set-TextValue "B47" "EnergySwap"
Set-TextValue "C47" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D47" "10.20"
Set-TextValue "E47" "  +0.09%  "
Set-TextValue "B48" "PaxDollar"
Set-TextValue "C48" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D48" "1.001"
Set-TextValue "E48" "  +0.23%  "
Set-TextValue "D49" "1.622"
Set-TextValue "E49" "  +1.51%  "
Set-TextValue "E50" "  +0.55%  "
Set-TextValue "D51" "64.16"
Set-TextValue "E51" "  +0.07%  "
